# iron_native / Excel COM-interop script reproducing the "Add working set of sequences" commit.
#
# The sheet holds a word/image/category cue table (rows 2-33, columns B/C/D,
# with A holding a 0-based running index). The edit swaps out the previous
# car/dog image+word sequence for a new house/dog sequence. Column A (index)
# and the header row (B1:D1) are unchanged; only B2:D33 get new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: image path for each cue, rows 2-33
$ws.Range("B2").Value = "house/house021.jpg"
$ws.Range("B3").Value = "dog/dog001.jpg"
$ws.Range("B4").Value = "dog/dog002.jpg"
$ws.Range("B5").Value = "dog/dog021.jpg"
$ws.Range("B6").Value = "dog/dog022.jpg"
$ws.Range("B7").Value = "house/house023.jpg"
$ws.Range("B8").Value = "dog/dog015.jpg"
$ws.Range("B9").Value = "dog/dog006.jpg"
$ws.Range("B10").Value = "house/house016.jpg"
$ws.Range("B11").Value = "dog/dog016.jpg"
$ws.Range("B12").Value = "dog/dog012.jpg"
$ws.Range("B13").Value = "house/house005.jpg"
$ws.Range("B14").Value = "dog/dog009.jpg"
$ws.Range("B15").Value = "house/house025.jpg"
$ws.Range("B16").Value = "dog/dog004.jpg"
$ws.Range("B17").Value = "house/house020.jpg"
$ws.Range("B18").Value = "house/house031.jpg"
$ws.Range("B19").Value = "house/house003.jpg"
$ws.Range("B20").Value = "house/house009.jpg"
$ws.Range("B21").Value = "house/house006.jpg"
$ws.Range("B22").Value = "dog/dog028.jpg"
$ws.Range("B23").Value = "dog/dog011.jpg"
$ws.Range("B24").Value = "dog/dog003.jpg"
$ws.Range("B25").Value = "house/house018.jpg"
$ws.Range("B26").Value = "house/house001.jpg"
$ws.Range("B27").Value = "dog/dog020.jpg"
$ws.Range("B28").Value = "dog/dog000.jpg"
$ws.Range("B29").Value = "house/house026.jpg"
$ws.Range("B30").Value = "house/house028.jpg"
$ws.Range("B31").Value = "dog/dog005.jpg"
$ws.Range("B32").Value = "house/house030.jpg"
$ws.Range("B33").Value = "house/house013.jpg"

# Column C: German verb/word for each cue, rows 2-33
$ws.Range("C2").Value = "sieben"
$ws.Range("C3").Value = "stärken"
$ws.Range("C4").Value = "wiegen"
$ws.Range("C5").Value = "pflegen"
$ws.Range("C6").Value = "krachen"
$ws.Range("C7").Value = "enden"
$ws.Range("C8").Value = "schicken"
$ws.Range("C9").Value = "opfern"
$ws.Range("C10").Value = "rasen"
$ws.Range("C11").Value = "gelten"
$ws.Range("C12").Value = "laufen"
$ws.Range("C13").Value = "gründen"
$ws.Range("C14").Value = "lehnen"
$ws.Range("C15").Value = "spielen"
$ws.Range("C16").Value = "bleiben"
$ws.Range("C17").Value = "kaufen"
$ws.Range("C18").Value = "töten"
$ws.Range("C19").Value = "stechen"
$ws.Range("C20").Value = "liefern"
$ws.Range("C21").Value = "füllen"
$ws.Range("C22").Value = "raten"
$ws.Range("C23").Value = "kehren"
$ws.Range("C24").Value = "loben"
$ws.Range("C25").Value = "drehen"
$ws.Range("C26").Value = "haken"
$ws.Range("C27").Value = "formen"
$ws.Range("C28").Value = "bitten"
$ws.Range("C29").Value = "währen"
$ws.Range("C30").Value = "hoffen"
$ws.Range("C31").Value = "tauschen"
$ws.Range("C32").Value = "schmecken"
$ws.Range("C33").Value = "runden"

# Column D: category label for each cue (matches the image folder), rows 2-33
$ws.Range("D2").Value = "house"
$ws.Range("D3").Value = "dog"
$ws.Range("D4").Value = "dog"
$ws.Range("D5").Value = "dog"
$ws.Range("D6").Value = "dog"
$ws.Range("D7").Value = "house"
$ws.Range("D8").Value = "dog"
$ws.Range("D9").Value = "dog"
$ws.Range("D10").Value = "house"
$ws.Range("D11").Value = "dog"
$ws.Range("D12").Value = "dog"
$ws.Range("D13").Value = "house"
$ws.Range("D14").Value = "dog"
$ws.Range("D15").Value = "house"
$ws.Range("D16").Value = "dog"
$ws.Range("D17").Value = "house"
$ws.Range("D18").Value = "house"
$ws.Range("D19").Value = "house"
$ws.Range("D20").Value = "house"
$ws.Range("D21").Value = "house"
$ws.Range("D22").Value = "dog"
$ws.Range("D23").Value = "dog"
$ws.Range("D24").Value = "dog"
$ws.Range("D25").Value = "house"
$ws.Range("D26").Value = "house"
$ws.Range("D27").Value = "dog"
$ws.Range("D28").Value = "dog"
$ws.Range("D29").Value = "house"
$ws.Range("D30").Value = "house"
$ws.Range("D31").Value = "dog"
$ws.Range("D32").Value = "house"
$ws.Range("D33").Value = "house"
